$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CERFileUpload")

# Clear the "Program" values (previously "-") in column I for rows 2-5,
# as part of setting up the batch definition / Events and Alerts page.
$ws.Range("I2").Value = $null
$ws.Range("I3").Value = $null
$ws.Range("I4").Value = $null
$ws.Range("I5").Value = $null
